$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.365.76'
$ws.Range('E2').Value = '  -1.60%  '
$ws.Range('D3').Value = '2.429.37'
$ws.Range('E3').Value = '  -1.56%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.28'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.20%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  -2.00%  '
$ws.Range('D9').Value = '2.425.62'
$ws.Range('E9').Value = '  -1.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.107'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.95%  '
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.21'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.351'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.56'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.82%  '
$ws.Range('E15').Value = '  -5.58%  '
$ws.Range('D16').Value = '2.873.35'
$ws.Range('E16').Value = '  -1.59%  '
$ws.Range('D17').Value = '62.278.20'
$ws.Range('E17').Value = '  -1.67%  '
$ws.Range('D18').Value = '2.436.97'
$ws.Range('E18').Value = '  -1.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.03'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.12'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.41%  '
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.09'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '620.23'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.93%  '
$ws.Range('D28').Value = '0.0₃0960'
$ws.Range('E28').Value = '  -8.67%  '
$ws.Range('E29').Value = '  -1.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.988'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.10%  '
$ws.Range('E31').Value = '  -5.11%  '
$ws.Range('E32').Value = '  -4.35%  '
$ws.Range('E33').Value = '  -2.18%  '
$ws.Range('E34').Value = '  -7.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.03'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.69%  '
$ws.Range('E36').Value = '  +0.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.44'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.79%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.374'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.51'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '147.01'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.21'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.27%  '
$ws.Range('E42').Value = '  -6.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.40'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.28%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.46'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '144.74'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.55%  '
$ws.Range('E47').Value = '  -1.91%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.07'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.31%  '
$ws.Range('E49').Value = '  -4.80%  '
$ws.Range('E50').Value = '  -2.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0228'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.15%  '
